$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / locator fields -------------------------------------------------
# Sales Order No.
$ws.Range("C14").Value = "sDB204-2311001"
# Target Date (weekly range display string)
$ws.Range("C19").Value = "27 Nov 2023 - 03 Dec 2023"
# Column header date range (By Invoice Date bucket)
$ws.Range("N22").Value = "04 Dec ~ 10 Dec"
# Purchase Order No. (repeated across the three line-item rows)
$ws.Range("E24").Value = "pDB204-2311001"
$ws.Range("E25").Value = "pDB204-2311001"
$ws.Range("E26").Value = "pDB204-2311001"

# --- Dates -------------------------------------------------------------
# Order Date
$ws.Range("C18").Value = 45250
# Row 23 forecast/plan dates
$ws.Range("P23").Value = 45252
$ws.Range("S23").Value = 45301
$ws.Range("T23").Value = 45343
$ws.Range("U23").Value = 45264
$ws.Range("V23").Value = 45265

# --- Quantity figures (InTransit / Receiver Inbounded / Estimated Inbound) --
# Row 24
$ws.Range("Q24").Value = 1620
$ws.Range("R24").Value = 0
$ws.Range("V24").Value = 1620
# Row 25
$ws.Range("Q25").Value = 1620
$ws.Range("R25").Value = 0
$ws.Range("V25").Value = 1620
# Row 26
$ws.Range("Q26").Value = 600
$ws.Range("R26").Value = 0
$ws.Range("V26").Value = 800
